{"js": "// Update the date line and all 25 \"three-digit \u00f7 one-digit\" problems in the table.\n// The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16) holds text.\n\nconst body = context.document.body;\n\n// 1) Update the date/weekday title paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text === \"2024-11-04 Monday\") {\n  titlePara.getRange().insertText(\"2024-11-05 Tuesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the division problems inside the table, cell by cell (row, col) so that\n//    identical old/new values across different cells never get cross-matched.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  \"0,0\": \"926\u00f75=\",\n  \"0,1\": \"941\u00f79=\",\n  \"0,2\": \"177\u00f76=\",\n  \"0,3\": \"948\u00f77=\",\n  \"0,4\": \"145\u00f78=\",\n  \"4,0\": \"647\u00f76=\",\n  \"4,1\": \"113\u00f79=\",\n  \"4,2\": \"695\u00f77=\",\n  \"4,3\": \"816\u00f74=\",\n  \"4,4\": \"517\u00f77=\",\n  \"8,0\": \"223\u00f79=\",\n  \"8,1\": \"564\u00f72=\",\n  \"8,2\": \"877\u00f73=\",\n  \"8,3\": \"262\u00f79=\",\n  \"8,4\": \"579\u00f77=\",\n  \"12,0\": \"125\u00f79=\",\n  \"12,1\": \"388\u00f79=\",\n  \"12,2\": \"889\u00f78=\",\n  \"12,3\": \"219\u00f78=\",\n  \"12,4\": \"713\u00f72=\",\n  \"16,0\": \"827\u00f74=\",\n  \"16,1\": \"324\u00f74=\",\n  \"16,2\": \"637\u00f78=\",\n  \"16,3\": \"432\u00f73=\",\n  \"16,4\": \"716\u00f72=\",\n};\n\nfor (const key of Object.keys(newValues)) {\n  const [row, col] = key.split(\",\").map(Number);\n  const cell = table.getCell(row, col);\n  const cellRange = cell.body.getRange();\n  cellRange.insertText(newValues[key], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date/weekday title (first paragraph in the document body).\n$titlePara = $d.Paragraphs(1)\nif ($titlePara.Range.Text.TrimEnd([char]13, [char]7) -eq \"2024-11-04 Monday\") {\n    $titlePara.Range.Text = \"2024-11-05 Tuesday\"\n}\n\n# 2) Update the division problems inside the table, cell by cell (row, col) so\n#    that identical old/new values across different cells never cross-match,\n#    which a blanket Find/Replace could do (e.g. 642\u00f75= -> 877\u00f73= and the\n#    pre-existing 877\u00f73= -> 637\u00f78= collide on the string \"877\u00f73=\").\n$t = $d.Tables(1)\n\n$newValues = @{\n    \"1,1\"  = \"926\u00f75=\";  \"1,2\"  = \"941\u00f79=\";  \"1,3\"  = \"177\u00f76=\";  \"1,4\"  = \"948\u00f77=\";  \"1,5\"  = \"145\u00f78=\";\n    \"5,1\"  = \"647\u00f76=\";  \"5,2\"  = \"113\u00f79=\";  \"5,3\"  = \"695\u00f77=\";  \"5,4\"  = \"816\u00f74=\";  \"5,5\"  = \"517\u00f77=\";\n    \"9,1\"  = \"223\u00f79=\";  \"9,2\"  = \"564\u00f72=\";  \"9,3\"  = \"877\u00f73=\";  \"9,4\"  = \"262\u00f79=\";  \"9,5\"  = \"579\u00f77=\";\n    \"13,1\" = \"125\u00f79=\";  \"13,2\" = \"388\u00f79=\";  \"13,3\" = \"889\u00f78=\";  \"13,4\" = \"219\u00f78=\";  \"13,5\" = \"713\u00f72=\";\n    \"17,1\" = \"827\u00f74=\";  \"17,2\" = \"324\u00f74=\";  \"17,3\" = \"637\u00f78=\";  \"17,4\" = \"432\u00f73=\";  \"17,5\" = \"716\u00f72=\";\n}\n\nforeach ($key in $newValues.Keys) {\n    $parts = $key.Split(\",\")\n    $row = [int]$parts[0]\n    $col = [int]$parts[1]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $newValues[$key]\n}\n"}
